{"js": "// Update the two-digit-divided-by-one-digit worksheet numbers.\n// The worksheet is a single 20x5 table where every 4th row (0, 4, 8, 12, 16)\n// holds the five division expressions for that block; the rows between are\n// blank spacer rows. We replace the text of each filled cell, in row-major\n// document order, with its new value.\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"isNullObject\");\nawait context.sync();\n\nif (table.isNullObject) {\n  throw new Error(\"Expected a table in the document body.\");\n}\n\n// [rowIndex, [old=>new for the 5 cells in that row]]\nconst rowUpdates = [\n  [0, [\"38\u00f72=\", \"56\u00f79=\", \"94\u00f76=\", \"52\u00f76=\", \"84\u00f75=\"]],\n  [4, [\"29\u00f72=\", \"95\u00f72=\", \"84\u00f72=\", \"45\u00f76=\", \"81\u00f75=\"]],\n  [8, [\"38\u00f76=\", \"18\u00f78=\", \"12\u00f78=\", \"82\u00f72=\", \"87\u00f74=\"]],\n  [12, [\"71\u00f79=\", \"18\u00f75=\", \"36\u00f79=\", \"10\u00f73=\", \"75\u00f73=\"]],\n  [16, [\"35\u00f74=\", \"18\u00f79=\", \"10\u00f77=\", \"52\u00f74=\", \"31\u00f73=\"]],\n];\n\nfor (const [rowIndex, newValues] of rowUpdates) {\n  for (let colIndex = 0; colIndex < newValues.length; colIndex++) {\n    const cell = table.getCellOrNullObject(rowIndex, colIndex);\n    cell.value = newValues[colIndex];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the two-digit-divided-by-one-digit worksheet numbers.\n# The worksheet is a single 20x5 table where every 4th row (Word's 1-based\n# rows 1, 5, 9, 13, 17) holds the five division expressions for that block;\n# the rows between are blank spacer rows. We overwrite the text of each\n# filled cell, in row-major document order, with its new value.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowUpdates = @(\n    @(1,  @(\"38\u00f72=\", \"56\u00f79=\", \"94\u00f76=\", \"52\u00f76=\", \"84\u00f75=\")),\n    @(5,  @(\"29\u00f72=\", \"95\u00f72=\", \"84\u00f72=\", \"45\u00f76=\", \"81\u00f75=\")),\n    @(9,  @(\"38\u00f76=\", \"18\u00f78=\", \"12\u00f78=\", \"82\u00f72=\", \"87\u00f74=\")),\n    @(13, @(\"71\u00f79=\", \"18\u00f75=\", \"36\u00f79=\", \"10\u00f73=\", \"75\u00f73=\")),\n    @(17, @(\"35\u00f74=\", \"18\u00f79=\", \"10\u00f77=\", \"52\u00f74=\", \"31\u00f73=\"))\n)\n\nforeach ($rowUpdate in $rowUpdates) {\n    $rowIndex = $rowUpdate[0]\n    $values = $rowUpdate[1]\n    for ($col = 1; $col -le $values.Count; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
